# "made hsa ready to run in cluster"
# Add a new experiment-configuration row (hsa_exp_index = 8) to Sheet1 that
# points at the cluster-mounted sample data path instead of the Windows
# workstation path used by the earlier rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 9

# Columns (row 1 headers):
# A hsa_exp_index | B data_type | C sub_data_type | D with_texture
# E exp_data_path | F file_ending | G age | H sex | I crop
# J crop_percentage | K landmark_placement | L manual_landmarks
# M export_landmarks | N calculate_hsa | O verbose
$ws.Cells.Item($newRow, 1).Value = 8
$ws.Cells.Item($newRow, 2).Value = "synthetic"
$ws.Cells.Item($newRow, 3).Value = "original"
$ws.Cells.Item($newRow, 4).Value = $true
$ws.Cells.Item($newRow, 5).Value = "/data/scratch/r092382/synthetic_data/synthetic_data_original_textured_unclipped_vtp_paraview_sample"
$ws.Cells.Item($newRow, 6).Value = ".vtp"
$ws.Cells.Item($newRow, 7).Value = 200
$ws.Cells.Item($newRow, 8).Value = "M"
$ws.Cells.Item($newRow, 9).Value = $false
$ws.Cells.Item($newRow, 10).Value = 0
$ws.Cells.Item($newRow, 11).Value = "automatic"
# Column L (manual_landmarks) intentionally left blank - placement is automatic.
$ws.Cells.Item($newRow, 13).Value = $false
$ws.Cells.Item($newRow, 14).Value = $true
$ws.Cells.Item($newRow, 15).Value = $true

# Move the selection to the freshly-added row, matching the author's final
# on-screen state after entering the new data.
$ws.Range("A9").Select()
